$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to stay/become plain text (matches the source data,
    # which stores these as inline strings) even when the new value looks
    # like a number (e.g. "362.00"); resetting the Style afterwards avoids
    # leaving a stray NumberFormat/quote-prefix style behind.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = '69.481.27'
$ws.Range("E2").Value = '  +2.43%  '

$ws.Range("D3").Value = '2.514.25'
$ws.Range("E3").Value = '  +0.85%  '

$ws.Range("E4").Value = '  +0.01%  '

Set-TextValue $ws.Range("D5") '597.89'
$ws.Range("E5").Value = '  +1.85%  '

Set-TextValue $ws.Range("D6") '176.29'
$ws.Range("E6").Value = '  -0.48%  '

$ws.Range("E7").Value = '  -0.06%  '

Set-TextValue $ws.Range("D8") '0.519'
$ws.Range("E8").Value = '  +0.77%  '

$ws.Range("D9").Value = '2.513.99'
$ws.Range("E9").Value = '  +0.90%  '

$ws.Range("E10").Value = '  +11.19%  '

$ws.Range("E11").Value = '  -0.35%  '

$ws.Range("E12").Value = '  +0.70%  '

$ws.Range("D14").Value = '2.973.09'
$ws.Range("E14").Value = '  +0.83%  '

Set-TextValue $ws.Range("D15") '25.97'
$ws.Range("E15").Value = '  +1.20%  '

$ws.Range("D18").Value = '2.507.75'
$ws.Range("E18").Value = '  -2.19%  '

Set-TextValue $ws.Range("D19") '7.66'
$ws.Range("E19").Value = '  +2.27%  '

Set-TextValue $ws.Range("D20") '362.00'
$ws.Range("E20").Value = '  +3.08%  '

Set-TextValue $ws.Range("D21") '11.01'
$ws.Range("E21").Value = '  +0.40%  '

Set-TextValue $ws.Range("D22") '4.08'
$ws.Range("E22").Value = '  -0.99%  '

$ws.Range("E23").Value = '  -0.10%  '

Set-TextValue $ws.Range("D24") '70.45'
$ws.Range("E24").Value = '  -0.54%  '

Set-TextValue $ws.Range("D25") '4.23'
$ws.Range("E25").Value = '  -1.18%  '

Set-TextValue $ws.Range("D26") '9.07'
$ws.Range("E26").Value = '  -0.46%  '

Set-TextValue $ws.Range("D27") '1.70'
$ws.Range("E27").Value = '  -2.75%  '

$ws.Range("D28").Value = '2.631.52'
$ws.Range("E28").Value = '  +0.42%  '

$ws.Range("E29").Value = '  +0.12%  '

Set-TextValue $ws.Range("D30") '511.44'
$ws.Range("E30").Value = '  +1.14%  '

$ws.Range("D31").Value = '0.0₃0896'
$ws.Range("E31").Value = '  -1.07%  '

Set-TextValue $ws.Range("D32") '7.74'
$ws.Range("E32").Value = '  -1.07%  '

$ws.Range("E33").Value = '  -1.68%  '

$ws.Range("E34").Value = '  +0.63%  '

$ws.Range("E35").Value = '  +0.05%  '

Set-TextValue $ws.Range("D36") '161.83'
$ws.Range("E36").Value = '  -0.88%  '

$ws.Range("E37").Value = '  -1.84%  '

Set-TextValue $ws.Range("D38") '18.72'
$ws.Range("E38").Value = '  +2.13%  '

$ws.Range("E40").Value = '  -1.04%  '

$ws.Range("E42").Value = '  -1.27%  '

$ws.Range("E43").Value = '  -1.40%  '

$ws.Range("E44").Value = '  -2.56%  '

Set-TextValue $ws.Range("D45") '2.36'
$ws.Range("E45").Value = '  -2.53%  '

Set-TextValue $ws.Range("D46") '38.81'
$ws.Range("E46").Value = '  -0.50%  '

Set-TextValue $ws.Range("D47") '149.90'
$ws.Range("E47").Value = '  +3.63%  '

$ws.Range("E48").Value = '  +1.72%  '

$ws.Range("E49").Value = '  +0.23%  '

$ws.Range("E50").Value = '  -1.22%  '

Set-TextValue $ws.Range("D51") '0.0739'
$ws.Range("E51").Value = '  -0.55%  '

# Row 16 and 17 swapped places (WrappedBTC <-> ShibaInu) with updated data
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue $ws.Range("D16") '0.0000178'
$ws.Range("E16").Value = '  +3.82%  '

$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '69.422.56'
$ws.Range("E17").Value = '  +2.52%  '
